# Paper/excel figures.xlsx -- "General Model Stats" sheet update
# - Table 1A (rows 2-12): iMR533/iMM518 comparison columns are swapped and
#   the "iMR533" model is renamed to "iMR534" with refreshed numbers.
# - Table 1B (rows 19-29): a handful of the iMR533 model statistics are
#   refreshed (Protein Coding Genes, metabolite/reaction counts).
# - Both table captions (B12 / B29) get updated text.
# - The now-empty trailing column E (rows 3-11) is cleared.
# - Selection is moved to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Model Stats")

# ---------------------------------------------------------------------
# Table 1A: Methanococcus maripaludis S2 model comparison (B2:D12)
#   column C = iMM518, column D = iMR534 (was C = iMR533, D = iMM518)
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "iMM518"
$ws.Range("D3").Value = "iMR534"

$ws.Range("C4").Value = 518
$ws.Range("D4").Value = 534

$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 31

$ws.Range("C6").Value = "556/49"
$ws.Range("D6").Value = "650/52"

$ws.Range("C7").Value = 163
$ws.Range("D7").Value = 268

$ws.Range("C8").Value = 570
$ws.Range("D8").Value = 571

$ws.Range("C9").Value = 49
$ws.Range("D9").Value = 57

$ws.Range("C10").Value = 464
$ws.Range("D10").Value = 570

$ws.Range("C11").Value = 75
$ws.Range("D11").Value = 91

# Table 1A caption: "roughly 100 more" -> "approximately 100 more"
$ws.Range("B12").Value = 'Table 1A. A comparison between iMR533 and iMM518 indicates that our model covers slightly more of the genome, including over 100 more gene-associated reactions. Both models include approximately the same number of reactions, but our model has approximately 100 more internal metabolites and dead end metabolites. Though this represent the portion of metabolism that cannot carry flux, all of our model''s dead end metabolites are part of gene-associated reactions and thus represent promising avenues for future model expansion. '

# The old trailing column E (rows 3-11) is no longer used.
$ws.Range("E3:E11").Clear()

# ---------------------------------------------------------------------
# Table 1B: Methanococcus maripaludis S2 model statistics (B19:C29)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 534
$ws.Range("C23").Value = "650/52"
$ws.Range("C24").Value = 268
$ws.Range("C25").Value = 571
$ws.Range("C27").Value = 220
$ws.Range("C28").Value = 570

# Table 1B caption: trimmed to a single sentence
$ws.Range("B29").Value = 'Table 1B. Some basic statistics for the iMR533 model.'

# ---------------------------------------------------------------------
# Final selection, matching the author's cursor position on save
# ---------------------------------------------------------------------
$ws.Range("B2:D17").Select()
